# Auto-generated edit script: update 2024 (column K / partial I,J) crime-count
# figures across the Citywide Totals, By Neighborhood, and per-neighborhood sheets
# to reflect the data refresh for 2024-06-20.

$wb = $excel.ActiveWorkbook

$changes = @{
    'Citywide Totals' = @(@{ Cell = "K2"; Value = 3597 }, @{ Cell = "K3"; Value = 3625 }, @{ Cell = "I4"; Value = 1795 }, @{ Cell = "J4"; Value = 1819 }, @{ Cell = "K4"; Value = 748 }, @{ Cell = "K5"; Value = 236 }, @{ Cell = "K6"; Value = 4200 }, @{ Cell = "I7"; Value = 26249 }, @{ Cell = "J7"; Value = 29287 }, @{ Cell = "K7"; Value = 12406 })
    'By Neighborhood' = @(@{ Cell = "K7"; Value = 354 }, @{ Cell = "K8"; Value = 833 }, @{ Cell = "K10"; Value = 68 }, @{ Cell = "K14"; Value = 67 }, @{ Cell = "K15"; Value = 127 }, @{ Cell = "K19"; Value = 380 }, @{ Cell = "K20"; Value = 281 }, @{ Cell = "K23"; Value = 123 }, @{ Cell = "K27"; Value = 124 }, @{ Cell = "K29"; Value = 654 }, @{ Cell = "K31"; Value = 138 }, @{ Cell = "K33"; Value = 509 }, @{ Cell = "K36"; Value = 155 }, @{ Cell = "K37"; Value = 419 }, @{ Cell = "K42"; Value = 439 }, @{ Cell = "K43"; Value = 113 }, @{ Cell = "J44"; Value = 232 }, @{ Cell = "K47"; Value = 69 }, @{ Cell = "K51"; Value = 146 }, @{ Cell = "K52"; Value = 340 }, @{ Cell = "K53"; Value = 167 }, @{ Cell = "K54"; Value = 244 }, @{ Cell = "K63"; Value = 47 }, @{ Cell = "K65"; Value = 297 }, @{ Cell = "K67"; Value = 485 }, @{ Cell = "K68"; Value = 27 }, @{ Cell = "K71"; Value = 38 }, @{ Cell = "K72"; Value = 61 }, @{ Cell = "K73"; Value = 114 }, @{ Cell = "K76"; Value = 183 }, @{ Cell = "K77"; Value = 87 }, @{ Cell = "K79"; Value = 313 }, @{ Cell = "J83"; Value = 593 }, @{ Cell = "K83"; Value = 264 }, @{ Cell = "K84"; Value = 87 }, @{ Cell = "K85"; Value = 566 }, @{ Cell = "I87"; Value = 71 }, @{ Cell = "K88"; Value = 144 }, @{ Cell = "K89"; Value = 171 }, @{ Cell = "K94"; Value = 149 }, @{ Cell = "K95"; Value = 207 }, @{ Cell = "K98"; Value = 63 }, @{ Cell = "K99"; Value = 214 }, @{ Cell = "I101"; Value = 26249 }, @{ Cell = "J101"; Value = 29287 }, @{ Cell = "K101"; Value = 12406 })
    'Logan Square' = @(@{ Cell = "K6"; Value = 85 }, @{ Cell = "K7"; Value = 167 })
    'Austin' = @(@{ Cell = "K2"; Value = 242 }, @{ Cell = "K4"; Value = 47 }, @{ Cell = "K6"; Value = 274 }, @{ Cell = "K7"; Value = 833 })
    'South Chicago' = @(@{ Cell = "K2"; Value = 95 }, @{ Cell = "K3"; Value = 88 }, @{ Cell = "J4"; Value = 25 }, @{ Cell = "J7"; Value = 593 }, @{ Cell = "K7"; Value = 264 })
    'Garfield Park' = @(@{ Cell = "K3"; Value = 193 }, @{ Cell = "K6"; Value = 145 }, @{ Cell = "K7"; Value = 509 })
    'West Pullman' = @(@{ Cell = "K3"; Value = 74 }, @{ Cell = "K7"; Value = 207 })
    'Grand Crossing' = @(@{ Cell = "K2"; Value = 111 }, @{ Cell = "K3"; Value = 145 }, @{ Cell = "K4"; Value = 22 }, @{ Cell = "K7"; Value = 419 })
    'New City' = @(@{ Cell = "K6"; Value = 119 }, @{ Cell = "K7"; Value = 297 })
    'Woodlawn' = @(@{ Cell = "K2"; Value = 60 }, @{ Cell = "K3"; Value = 83 }, @{ Cell = "K7"; Value = 214 })
    'Gage Park' = @(@{ Cell = "K6"; Value = 50 }, @{ Cell = "K7"; Value = 138 })
    'North Lawndale' = @(@{ Cell = "K3"; Value = 160 }, @{ Cell = "K7"; Value = 485 })
    'South Deering' = @(@{ Cell = "K6"; Value = 21 }, @{ Cell = "K7"; Value = 87 })
    'Loop' = @(@{ Cell = "K3"; Value = 71 }, @{ Cell = "K7"; Value = 244 })
    'Englewood' = @(@{ Cell = "K2"; Value = 187 }, @{ Cell = "K3"; Value = 225 }, @{ Cell = "K5"; Value = 16 }, @{ Cell = "K7"; Value = 654 })
    'Chatham' = @(@{ Cell = "K3"; Value = 108 }, @{ Cell = "K7"; Value = 380 })
    'Irving Park' = @(@{ Cell = "J4"; Value = 14 }, @{ Cell = "J7"; Value = 232 })
    'River North' = @(@{ Cell = "K6"; Value = 106 }, @{ Cell = "K7"; Value = 183 })
    'Bridgeport' = @(@{ Cell = "K3"; Value = 12 }, @{ Cell = "K6"; Value = 22 }, @{ Cell = "K7"; Value = 67 })
    'Humboldt Park' = @(@{ Cell = "K3"; Value = 143 }, @{ Cell = "K7"; Value = 439 })
    'Avondale' = @(@{ Cell = "K6"; Value = 31 }, @{ Cell = "K7"; Value = 68 })
    'Rogers Park' = @(@{ Cell = "K2"; Value = 44 }, @{ Cell = "K6"; Value = 58 })
    'Douglas' = @(@{ Cell = "K5"; Value = 2 }, @{ Cell = "K7"; Value = 123 })
    'Roseland' = @(@{ Cell = "K2"; Value = 109 }, @{ Cell = "K7"; Value = 313 })
    'Chicago Lawn' = @(@{ Cell = "K2"; Value = 96 }, @{ Cell = "K3"; Value = 83 }, @{ Cell = "K7"; Value = 281 })
    'Grand Boulevard' = @(@{ Cell = "K3"; Value = 45 }, @{ Cell = "K6"; Value = 33 }, @{ Cell = "K7"; Value = 155 })
    'Auburn Gresham' = @(@{ Cell = "K2"; Value = 131 }, @{ Cell = "K7"; Value = 354 })
    'West Loop' = @(@{ Cell = "K2"; Value = 42 }, @{ Cell = "K7"; Value = 149 })
    'Kenwood' = @(@{ Cell = "K2"; Value = 21 }, @{ Cell = "K7"; Value = 69 })
    'Brighton Park' = @(@{ Cell = "K2"; Value = 45 }, @{ Cell = "K3"; Value = 32 }, @{ Cell = "K7"; Value = 127 })
    'Wicker Park' = @(@{ Cell = "K3"; Value = 10 }, @{ Cell = "K7"; Value = 63 })
    'Portage Park' = @(@{ Cell = "K2"; Value = 34 }, @{ Cell = "K3"; Value = 27 }, @{ Cell = "K6"; Value = 46 }, @{ Cell = "K7"; Value = 114 })
    'United Center' = @(@{ Cell = "K6"; Value = 68 }, @{ Cell = "K7"; Value = 144 })
    'Uptown' = @(@{ Cell = "K2"; Value = 41 }, @{ Cell = "K3"; Value = 54 }, @{ Cell = "K6"; Value = 52 }, @{ Cell = "K7"; Value = 171 })
    'Edgewater' = @(@{ Cell = "K6"; Value = 47 }, @{ Cell = "K7"; Value = 124 })
    'Little Italy, UIC' = @(@{ Cell = "K3"; Value = 41 }, @{ Cell = "K7"; Value = 146 })
    'North Park' = @(@{ Cell = "K2"; Value = 9 }, @{ Cell = "K7"; Value = 27 })
    'Hyde Park' = @(@{ Cell = "K2"; Value = 23 }, @{ Cell = "K3"; Value = 30 }, @{ Cell = "K7"; Value = 113 })
    'South Shore' = @(@{ Cell = "K2"; Value = 200 }, @{ Cell = "K3"; Value = 191 }, @{ Cell = "K6"; Value = 131 }, @{ Cell = "K7"; Value = 566 })
    'Oakland' = @(@{ Cell = "K3"; Value = 11 }, @{ Cell = "K7"; Value = 38 })
    'Old Town' = @(@{ Cell = "K2"; Value = 7 }, @{ Cell = "K7"; Value = 61 })
    'Riverdale' = @(@{ Cell = "K3"; Value = 31 }, @{ Cell = "K7"; Value = 87 })
    'Little Village' = @(@{ Cell = "K2"; Value = 87 }, @{ Cell = "K7"; Value = 340 })
    'Ukrainian Village' = @(@{ Cell = "I4"; Value = 5 }, @{ Cell = "I7"; Value = 71 })
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($change in $changes[$sheetName]) {
        $ws.Range($change.Cell).Value = $change.Value
    }
}

Write-Output "Updated $($changes.Keys.Count) sheets"
